$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (周欣媛 / 212241811426): new article + new Materials & Methods note ---
$ws.Range("D4").Value = "Jones, Carl._2014_Science Fiction Studies_Stealing Kinship: Neuromancer and Artificial Intelligence"

$ws.Range("G4").Value = @"
materials
1. Primary Text Analysis
2. Secondary Literature
3. Cultural and Technological Context
methods
1. Literary Analysis
2. Interdisciplinary Approach
3. Comparative Analysis
4. Cultural Critique
"@

# --- Clear the stray "literary analysis" notes from the Appendix column ---
$ws.Range("K3").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("K6").ClearContents()

# --- View state: zoom in and move the selection/scroll position ---
$excel.ActiveWindow.Zoom = 40
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F4").Select()
